# Update parameter threshold table:
#  - alpha_distance_range : Min/Max refreshed
#  - beta_distance_range  : Min/Max refreshed
#  - ratio_threshold_range: Min/Max refreshed
#  - theta_threshold_range: entire row removed
#  - pie_threshold_range  : Min/Max refreshed (row shifts up to row 5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Min/Max values for the rows that remain (rows 2-4 keep their
# position; row 5 "theta_threshold_range" is removed below which shifts the
# former row 6 "pie_threshold_range" up to row 5).
$ws.Range("B2").Value = 5.4
$ws.Range("C2").Value = 11.7

$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 10.3

$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.2

# Remove the whole "theta_threshold_range" row (row 5). This shifts the
# "pie_threshold_range" row up from row 6 to row 5, shrinking the used range
# to A1:C5 and dropping the now-unused shared string / font.
$ws.Rows.Item(5).Delete()

# Set the new Min/Max for "pie_threshold_range" (now row 5).
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Match the saved selection/active cell recorded in the workbook.
[void]$ws.Range("E6").Select()
